# Apply "fingerprint demo" update to testprotocol.xlsx
#
# Summary of changes (per the target diff):
#  - sheet "protocoltestcasedetails": D19 (testcase18) flips N -> Y,
#    D21 (testcase20) flips Y -> N.
#  - The active/selected sheet moves from "protocol" to
#    "protocoltestcasedetails" (tabSelected / workbookView activeTab).
#  - Selection on "protocol" moves from B8 to B12.
#  - Selection on "protocoltestcasedetails" moves from C19 to E24,
#    and its view scrolls down (topLeftCell A10 -> A13).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("protocol")
$ws2 = $wb.Worksheets.Item("protocoltestcasedetails")

# --- Cell value changes on protocoltestcasedetails ---
# D19 (testcase18_parquet_dbtable_match_likeobject): N -> Y
$ws2.Range("D19").Value = "Y"
# D21 (testcase20_oracle_bigquery_match_manual): Y -> N
$ws2.Range("D21").Value = "N"

# --- Update selection left behind on the first sheet (no longer active) ---
$ws1.Activate()
$ws1.Range("B12").Select()

# --- Make protocoltestcasedetails the active sheet, with the new
#     selection/scroll position, matching the saved view state. ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("E24").Select()

Write-Host "Applied fingerprint-demo updates to testprotocol.xlsx"
